# Apply "Updated problem statement and time planning sheet" edit.
# Five cards that were marked as "Text Cards to implement:" (column C) have
# now been finished, so they move into the "Finished Cards:" list
# (column D), landing in new rows to keep column D's alphabetical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the newly finished cards to column D first (so the shared string
# entries stay alive / stable), then clear their old column C entries.
$ws.Range("D10").Value = "Library"
$ws.Range("D11").Value = "Throne Room"
$ws.Range("D23").Value = "Mine"
$ws.Range("D24").Value = "Adventurer"
$ws.Range("D25").Value = "Chapel"

$ws.Range("C9").ClearContents()
$ws.Range("C13").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("C25").ClearContents()

# Update the saved selection/active cell shown when the sheet is reopened.
$null = $ws.Range("C27").Select()
